$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: D14 changed 105 -> 92 (C14 = SUM(D14:E14) recomputes automatically)
$ws.Range("D14").Value = 92

# Row 17: F17, I17, J17 changed
$ws.Range("F17").Value = 24
$ws.Range("I17").Value = 0.7
$ws.Range("J17").Value = 0.15

# Row 20: D20, E20 changed
$ws.Range("D20").Value = 66
$ws.Range("E20").Value = 22

# Row 22: D22, E22 changed
$ws.Range("D22").Value = 75
$ws.Range("E22").Value = 14

# P22 text changed from placeholder note to citation
$ws.Range("P22").Value = "Pelletier 2021 In prep"

# Update sheet view: remove topLeftCell scroll position, change selection to M22
$ws.Range("M22").Select()
